$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Developer name
$ws.Range("C3").Value2 = "Om Patel"

# Bold the column headers for Preconditions / Method Inputs / Expected Result
$ws.Range("E6:G6").Font.Bold = $true

# Row 7 - __init__ : Attributes are set to input values.
$ws.Range("E7").Value2 = "None "
$ws.Range("F7").Value2 = "account_number = 2004, client_number = 2904, balance = 1000.00"
$ws.Range("G7").Value2 = "Attributes are correctly initialized."

# Row 8 - __init__ : Balance attribute set to 0 when non-numeric balance argument.
$ws.Range("E8").Value2 = "None "
$ws.Range("F8").Value2 = 'account_number=2004, client_number=2904, balance="not_a_number"'
$ws.Range("G8").Value2 = "Balance is set to 0.00"

# Row 9 - __init__ : ValueError when non-numeric account number
$ws.Range("E9").Value2 = "None "
$ws.Range("F9").Value2 = 'account_number="not_a_number", client_number=2904, balance=1000.00'
$ws.Range("G9").Value2 = 'Raises ValueError with message "Account number must be an integer."'

# Row 10 - __init__ : ValueError when non-numeric client number
$ws.Range("E10").Value2 = "None "
$ws.Range("F10").Value2 = 'account_number=2004, client_number="not_a_number", balance=1000.00'
$ws.Range("G10").Value2 = "Raises ValueError"

# Row 11 - account_number (getter) : returns account number attribute
$ws.Range("E11").Value2 = "Account is initialized"
$ws.Range("F11").Value2 = "None"
$ws.Range("G11").Value2 = "Return 2004"

# Row 12 - client_number (getter) : returns client number attribute
$ws.Range("E12").Value2 = "Account is initialized"
$ws.Range("F12").Value2 = "None"
$ws.Range("G12").Value2 = "Return 2904"

# Row 13 - balance (getter) : returns balance attribute
$ws.Range("E13").Value2 = "Account is initialized"
$ws.Range("F13").Value2 = "None"
$ws.Range("G13").Value2 = "Return 1000.00"

# Row 14 - update_balance : correctly updates balance attribute when positive amount is received.
$ws.Range("E14").Value2 = "Account is initialized with balance of 1000.00"
$ws.Range("F14").Value2 = "amount = 500.00"
$ws.Range("G14").Value2 = "Balance becomes 1500.00"

# Row 15 - update_balance : correctly updates balance attribute when negative amount is received.
$ws.Range("E15").Value2 = "Account is initialized with balance of 1000.00"
$ws.Range("F15").Value2 = "amount = - 200"
$ws.Range("G15").Value2 = "Balance becomes 800.00"

# Row 16 - update_balance : Balance attribute value remains unchanged when amount is non-numeric
$ws.Range("E16").Value2 = "Account is initialized with balance of 1000.00"
$ws.Range("F16").Value2 = 'amount = "not_a_number"'
$ws.Range("G16").Value2 = "Balance remains 1000.00"

# Row 17 - deposit : BankAccount object's balance is updated correctly when a valid amount is provided.
$ws.Range("E17").Value2 = "Account is initialized with balance of 1000.00"
$ws.Range("F17").Value2 = "amount = 150.00"
$ws.Range("G17").Value2 = "Balance becomes 1150.00"

# Row 18 - deposit : ValueError when negative amount is provided.
$ws.Range("E18").Value2 = "Account is initialized "
$ws.Range("F18").Value2 = "amount = - 50.00"
$ws.Range("G18").Value2 = 'Raises ValueError with message "Deposit amount must be positive."'

# Row 19 - withdraw : BankAccount object's balance is updated correctly when a valid amount is provided.
$ws.Range("E19").Value2 = "Account is initialized with balance of 1000.00"
$ws.Range("F19").Value2 = "amount = 100.00"
$ws.Range("G19").Value2 = "Balance becomes 900.00"

# Row 20 - withdraw : ValueError when negative amount is provided.
$ws.Range("E20").Value2 = "Account is initialized"
$ws.Range("F20").Value2 = "amount = - 50.00"
$ws.Range("G20").Value2 = 'Raises ValueError with message "Withdrawal amount must be positive."'

# Row 21 - withdraw : ValueError when amount exceeds balance.
$ws.Range("E21").Value2 = "Account is initialized with balance of 1000.00"
$ws.Range("F21").Value2 = "amount = 200.00"
$ws.Range("G21").Value2 = 'Raises ValueError with message "Withdrawal amount must not exceed the account balance."'

# Row 22 - __str__ : returns string in expected format.
$ws.Range("E22").Value2 = "Account is initialized with account_number=2004, client_number=2904, and balance=1000.00"
$ws.Range("F22").Value2 = "None"
$ws.Range("G22").Value2 = 'Returns "Account: 2004, Client: 2904, Balance: 1000.00"'
